# Regen save_data to use K (strikeouts) instead of Strike# in column G.
# Updates the K column (G2:G21) with the recalculated strikeout counts.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 4
    3  = 6
    4  = 6
    5  = 6
    6  = 6
    7  = 1
    8  = 1
    9  = 3
    10 = 2
    11 = 7
    12 = 2
    13 = 2
    14 = 2
    15 = 6
    16 = 4
    17 = 7
    18 = 1
    19 = 5
    20 = 1
    21 = 2
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
